$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize legacy float-typed numeric cells in row 2 to clean integers
$ws.Range("A2").Value = 1
$ws.Range("F2").Value = 12345
$ws.Range("J2").Value = 1

# New test user row (row 11) for quick login testing
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "test"
$ws.Range("C11").Value = "Test"
$ws.Range("D11").Value = "Test"
$ws.Range("E11").Value = "Teststraße 1"
$ws.Range("F11").Value = 77777
$ws.Range("G11").Value = "Teststadt"
$ws.Range("H11").Value = "test.test@example.com"
$ws.Hyperlinks.Add($ws.Range("H11"), "mailto:test.test@example.com")
$ws.Range("I11").Value = "test"
$ws.Range("I11").NumberFormat = $ws.Range("I2").NumberFormat
$ws.Range("I11").HorizontalAlignment = -4108
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = $false

$ws.Range("K18").Select() | Out-Null

Write-Host "Test user row added"
